{"js": "// Apply the worksheet's text updates: refresh the date stamp and swap in a\n// fresh batch of three-digit-by-one-digit division problems.\n// Each entry is unique in the document, so a scoped search-and-replace on\n// the exact original text is safe and order independent.\nconst replacements = [\n  [\"2025-06-26 Thursday\", \"2025-06-27 Friday\"],\n  [\"911\u00f75=\", \"540\u00f76=\"],\n  [\"513\u00f73=\", \"711\u00f73=\"],\n  [\"225\u00f76=\", \"918\u00f79=\"],\n  [\"356\u00f77=\", \"549\u00f79=\"],\n  [\"297\u00f76=\", \"407\u00f73=\"],\n  [\"751\u00f75=\", \"180\u00f77=\"],\n  [\"146\u00f74=\", \"891\u00f73=\"],\n  [\"178\u00f76=\", \"498\u00f76=\"],\n  [\"764\u00f74=\", \"687\u00f72=\"],\n  [\"427\u00f73=\", \"845\u00f77=\"],\n  [\"725\u00f79=\", \"591\u00f76=\"],\n  [\"974\u00f72=\", \"603\u00f74=\"],\n  [\"725\u00f74=\", \"222\u00f75=\"],\n  [\"913\u00f78=\", \"263\u00f74=\"],\n  [\"263\u00f74=\", \"571\u00f79=\"],\n  [\"947\u00f76=\", \"688\u00f79=\"],\n  [\"974\u00f76=\", \"899\u00f78=\"],\n  [\"898\u00f72=\", \"969\u00f78=\"],\n  [\"645\u00f75=\", \"424\u00f75=\"],\n  [\"563\u00f73=\", \"467\u00f75=\"],\n  [\"615\u00f77=\", \"381\u00f76=\"],\n  [\"797\u00f75=\", \"837\u00f79=\"],\n  [\"131\u00f72=\", \"343\u00f74=\"],\n  [\"759\u00f72=\", \"559\u00f72=\"],\n  [\"280\u00f73=\", \"759\u00f79=\"],\n];\n\nconst body = context.document.body;\n\n// Two-pass replace: some new values equal other cells' old values (e.g. the\n// cell that read \"913\u00f78=\" becomes \"263\u00f74=\", while the cell that already read\n// \"263\u00f74=\" becomes \"571\u00f79=\"). Writing directly in one pass would make the\n// second search also match the text the first pass just wrote. Route every\n// value through a unique placeholder first, then resolve placeholders to the\n// real final text in a second pass.\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText] = replacements[i];\n  const placeholder = \"\\u0001PLACEHOLDER\" + i + \"\\u0001\";\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(placeholder, \"Replace\");\n  }\n  await context.sync();\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const placeholder = \"\\u0001PLACEHOLDER\" + i + \"\\u0001\";\n  const results = body.search(placeholder, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the worksheet's text updates: refresh the date stamp and swap in a\n# fresh batch of three-digit-by-one-digit division problems.\n#\n# Two-pass replace: some new values equal other cells' old values (e.g. the\n# cell that read \"913\u00f78=\" becomes \"263\u00f74=\", while the cell that already read\n# \"263\u00f74=\" becomes \"571\u00f79=\"). Replacing directly in one pass would make the\n# second Find also match the text the first pass just wrote. Route every\n# value through a unique placeholder first, then resolve placeholders to the\n# real final text in a second pass.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-26 Thursday\", \"2025-06-27 Friday\"),\n    @(\"911\u00f75=\", \"540\u00f76=\"),\n    @(\"513\u00f73=\", \"711\u00f73=\"),\n    @(\"225\u00f76=\", \"918\u00f79=\"),\n    @(\"356\u00f77=\", \"549\u00f79=\"),\n    @(\"297\u00f76=\", \"407\u00f73=\"),\n    @(\"751\u00f75=\", \"180\u00f77=\"),\n    @(\"146\u00f74=\", \"891\u00f73=\"),\n    @(\"178\u00f76=\", \"498\u00f76=\"),\n    @(\"764\u00f74=\", \"687\u00f72=\"),\n    @(\"427\u00f73=\", \"845\u00f77=\"),\n    @(\"725\u00f79=\", \"591\u00f76=\"),\n    @(\"974\u00f72=\", \"603\u00f74=\"),\n    @(\"725\u00f74=\", \"222\u00f75=\"),\n    @(\"913\u00f78=\", \"263\u00f74=\"),\n    @(\"263\u00f74=\", \"571\u00f79=\"),\n    @(\"947\u00f76=\", \"688\u00f79=\"),\n    @(\"974\u00f76=\", \"899\u00f78=\"),\n    @(\"898\u00f72=\", \"969\u00f78=\"),\n    @(\"645\u00f75=\", \"424\u00f75=\"),\n    @(\"563\u00f73=\", \"467\u00f75=\"),\n    @(\"615\u00f77=\", \"381\u00f76=\"),\n    @(\"797\u00f75=\", \"837\u00f79=\"),\n    @(\"131\u00f72=\", \"343\u00f74=\"),\n    @(\"759\u00f72=\", \"559\u00f72=\"),\n    @(\"280\u00f73=\", \"759\u00f79=\")\n)\n\n# Pass 1: old text -> unique placeholder\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $oldText = $replacements[$i][0]\n    $placeholder = \"@@PLACEHOLDER\" + $i + \"@@\"\n\n    $find = $d.Content.Find\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2)\n}\n\n# Pass 2: placeholder -> final new text\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $newText = $replacements[$i][1]\n    $placeholder = \"@@PLACEHOLDER\" + $i + \"@@\"\n\n    $find = $d.Content.Find\n    $find.Execute($placeholder, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
